$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (RandomForestRegressor)
$ws.Range("B3").Value = 0.9929497140304089
$ws.Range("C3").Value = 0.9931427002564103
$ws.Range("D3").Value = 0.9929636102501272

# Row 4 (GradientBoostingRegressor)
$ws.Range("B4").Value = 0.9953808967668092
$ws.Range("C4").Value = 0.9953808650550096
$ws.Range("D4").Value = 0.995380909138542

# Row 5 (AdaBoostRegressor)
$ws.Range("B5").Value = 0.9860014569604251
$ws.Range("C5").Value = 0.9849701780848438
$ws.Range("D5").Value = 0.9854476329559695
